$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.031.12'
$ws.Range('E2').Value = '  -3.90%  '

$ws.Range('D3').Value = '2.496.99'
$ws.Range('E3').Value = '  -5.08%  '

$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').Value = '544.55'
$ws.Range('E5').Value = '  -1.44%  '

$ws.Range('D6').Value = '147.38'
$ws.Range('E6').Value = '  -4.66%  '

$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.24%  '

$ws.Range('D8').Value = '0.583'

$ws.Range('D9').Value = '2.524.96'
$ws.Range('E9').Value = '  -4.13%  '

$ws.Range('E10').Value = '  -3.21%  '

$ws.Range('E11').Value = '  -0.87%  '

$ws.Range('D12').Value = '5.52'
$ws.Range('E12').Value = '  +2.21%  '

$ws.Range('D13').Value = '0.357'
$ws.Range('E13').Value = '  -1.74%  '

$ws.Range('D14').Value = '2.953.53'
$ws.Range('E14').Value = '  -4.65%  '

$ws.Range('E15').Value = '  -4.27%  '

$ws.Range('D16').Value = '60.037.28'
$ws.Range('E16').Value = '  -3.74%  '

$ws.Range('E17').Value = '  -2.48%  '

$ws.Range('D18').Value = '2.506.57'
$ws.Range('E18').Value = '  -4.76%  '

$ws.Range('D19').Value = '11.48'
$ws.Range('E19').Value = '  -1.75%  '

$ws.Range('D20').Value = '4.38'
$ws.Range('E20').Value = '  -3.42%  '

$ws.Range('D21').Value = '328.26'
$ws.Range('E21').Value = '  -3.50%  '

$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.01%  '

$ws.Range('D23').Value = '5.84'
$ws.Range('E23').Value = '  -4.37%  '

$ws.Range('D24').Value = '61.71'
$ws.Range('E24').Value = '  -1.74%  '

$ws.Range('D25').Value = '0.450'
$ws.Range('E25').Value = '  -10.29%  '

$ws.Range('E26').Value = '  +1.32%  '

$ws.Range('D27').Value = '0.163'
$ws.Range('E27').Value = '  -3.40%  '

$ws.Range('D28').Value = '7.87'
$ws.Range('E28').Value = '  -2.20%  '

$ws.Range('D29').Value = '1.31'
$ws.Range('E29').Value = '  -1.83%  '

$ws.Range('D30').Value = '0.0₃0798'
$ws.Range('E30').Value = '  -4.63%  '

$ws.Range('D31').Value = '6.91'
$ws.Range('E31').Value = '  -2.58%  '

$ws.Range('E32').Value = '  -3.43%  '

$ws.Range('D33').Value = '0.997'
$ws.Range('E33').Value = '  -0.17%  '

$ws.Range('D34').Value = '158.87'
$ws.Range('E34').Value = '  -1.35%  '

$ws.Range('D35').Value = '1.43'
$ws.Range('E35').Value = '  +0.27%  '

$ws.Range('D36').Value = '18.97'
$ws.Range('E36').Value = '  -1.39%  '

$ws.Range('E37').Value = '  -4.23%  '

$ws.Range('D38').Value = '1.75'
$ws.Range('E38').Value = '  +0.65%  '

$ws.Range('E39').Value = '  -0.57%  '

$ws.Range('D40').Value = '315.35'
$ws.Range('E40').Value = '  -6.18%  '

$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').Value = '36.81'
$ws.Range('E41').Value = '  -3.00%  '

$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = '3.80'
$ws.Range('E42').Value = '  -3.33%  '

$ws.Range('D43').Value = '0.842'
$ws.Range('E43').Value = '  -7.42%  '

$ws.Range('D44').Value = '0.995'
$ws.Range('E44').Value = '  -0.33%  '

$ws.Range('D45').Value = '0.606'
$ws.Range('E45').Value = '  -1.07%  '

$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').Value = '10.78'
$ws.Range('E46').Value = '  -1.75%  '

$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '127.44'
$ws.Range('E47').Value = '  -0.07%  '

$ws.Range('E48').Value = '  -2.21%  '

$ws.Range('D49').Value = '0.0947'
$ws.Range('E49').Value = '  -1.62%  '

$ws.Range('E50').Value = '  -2.18%  '

$ws.Range('D51').Value = '18.75'
$ws.Range('E51').Value = '  -5.27%  '
